$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("MODERN SPORTS CARS")
$ws.Range("A2").Value = "2017 Abarth 124 SPIDER "
$ws.Range("A3").Value = "2014 Alfa Romeo 4C "
$ws.Range("A4").Value = "2017 Alpine A110 "

$ws = $wb.Worksheets.Item("HOT HATCH")
$ws.Range("A2").Value = "2016 Abarth 695 BIPOSTO "

$ws = $wb.Worksheets.Item("CLASSIC RALLY")
$ws.Range("A2").Value = "1980 Abarth FIAT 131 "

$ws = $wb.Worksheets.Item("CULT CARS")
$ws.Range("A2").Value = "1968 Abarth 595 ESSEESSE "

$ws = $wb.Worksheets.Item("MODERN SUPERCARS")
$ws.Range("A2").Value = "2017 Acura NSX "

$ws = $wb.Worksheets.Item("RETRO HOT HATCH")
$ws.Range("A2").Value = "2002 Acura RSX TYPE-S "
$ws.Range("A3").Value = "2001 Acura INTEGRA TYPE-R "

$ws = $wb.Worksheets.Item("SPORTS UTILITY HEROES")
$ws.Range("A2").Value = "2018 Alfa Romeo STELVIO QUADRIFOGLIO "

$ws = $wb.Worksheets.Item("SUPER SALOONS")
$ws.Range("A2").Value = "2017 Alfa Romeo GIULIA QUADRIFOGLIO "
$ws.Range("A3").Value = "2016 Alfa Romeo GIULIA QUADRIFOGLIO FORZA EDITION "

$ws = $wb.Worksheets.Item("GT CARS")
$ws.Range("A2").Value = "2007 Alfa Romeo 8C COMPETIZIONE "
$ws.Range("A3").Value = "2007 Alfa Romeo 8C COMPETIZIONE FORZA EDITION "

$ws = $wb.Worksheets.Item("RETRO SALOONS")
$ws.Range("A2").Value = "1992 Alfa Romeo 155 Q4 "

$ws = $wb.Worksheets.Item("CLASSIC RACERS")
$ws.Range("A2").Value = "1968 Alfa Romeo 33 STRADALE "
$ws.Range("A3").Value = "1965 Alfa Romeo GIULIA TZ2 "

$ws = $wb.Worksheets.Item("RARE CLASSICS")
$ws.Range("A2").Value = "1965 Alfa Romeo GIULIA SPRINT GTA STRADALE "

$ws = $wb.Worksheets.Item("VINTAGE RACERS")
$ws.Range("A2").Value = "1934 Alfa Romeo P3 "

$ws = $wb.Worksheets.Item("OFFROADS BUGGIES")
$ws.Range("A2").Value = "2015 Alumi Craft CLASS 10 RACE CAR "
